# Update cryptos list with latest price/volume data from coinranking.com
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 15 and 16 swap rank: Chainlink moves up to rank 13 (row 15),
# WrappedEther moves down to rank 14 (row 16)
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.586"
$ws.Range("E15").Value = "  -3.30%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.551.37"
$ws.Range("E16").Value = "  +0.00%  "

# Updated price and volume figures for remaining rows
$ws.Range("D2").Value = "21.648.38"
$ws.Range("E2").Value = "  -1.96%  "
$ws.Range("D3").Value = "1.534.61"
$ws.Range("E3").Value = "  -1.52%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "288.41"
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3929"
$ws.Range("E7").Value = "  +1.21%  "
$ws.Range("E8").Value = "  -2.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.42"
$ws.Range("E9").Value = "  -1.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07176"
$ws.Range("E10").Value = "  -2.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.055"
$ws.Range("E11").Value = "  -6.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.673"
$ws.Range("E13").Value = "  -0.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.56"
$ws.Range("E14").Value = "  -4.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001094"
$ws.Range("E17").Value = "  -3.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06586"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "83.69"
$ws.Range("E19").Value = "  -1.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9999"
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.119"
$ws.Range("E21").Value = "  -4.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.46"
$ws.Range("E22").Value = "  -3.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.71"
$ws.Range("E23").Value = "  -6.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.357"
$ws.Range("E24").Value = "  +0.64%  "
$ws.Range("D25").Value = "21.658.01"
$ws.Range("E25").Value = "  -1.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.354"
$ws.Range("E26").Value = "  -8.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "149.37"
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("E28").Value = "  -3.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.851"
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("D30").Value = "1.712.79"
$ws.Range("E30").Value = "  -0.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.06"
$ws.Range("E31").Value = "  -3.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.079"
$ws.Range("E32").Value = "  +3.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9389"
$ws.Range("E33").Value = "  -15.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08142"
$ws.Range("E34").Value = "  -0.90%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.518"
$ws.Range("E35").Value = "  -8.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.168"
$ws.Range("E36").Value = "  -1.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06021"
$ws.Range("E37").Value = "  -4.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02220"
$ws.Range("E38").Value = "  -3.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.458"
$ws.Range("E39").Value = "  -14.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2023"
$ws.Range("E40").Value = "  -4.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.181"
$ws.Range("E41").Value = "  -3.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.91"
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9996"
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5766"
$ws.Range("E44").Value = "  -3.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.07"
$ws.Range("E45").Value = "  -3.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.707"
$ws.Range("E46").Value = "  -0.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5505"
$ws.Range("E47").Value = "  -4.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.163"
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("E49").Value = "  -2.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "115.98"
$ws.Range("E50").Value = "  -2.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06694"
$ws.Range("E51").Value = "  -3.04%  "
Write-Host "Cryptos list updated"
